{"js": "// Update due-date strings for the Technology & Service Inspiration Posts\n// and Reflection Posts lists.\nconst replacements = [\n  [\"September 04\", \"February 3\"],\n  [\"September 25\", \"February 24\"],\n  [\"October 30\", \"April 06\"],\n  [\"September 11\", \"February 10\"],\n  [\"October 02\", \"March 02\"],\n  [\"November 20\", \"April 27\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"September 04\"; New = \"February 3\" },\n    @{ Old = \"September 25\"; New = \"February 24\" },\n    @{ Old = \"October 30\";   New = \"April 06\" },\n    @{ Old = \"September 11\"; New = \"February 10\" },\n    @{ Old = \"October 02\";   New = \"March 02\" },\n    @{ Old = \"November 20\";  New = \"April 27\" }\n)\n\nforeach ($p in $d.Paragraphs) {\n    foreach ($r in $replacements) {\n        if ($p.Range.Text -eq ($r.Old + \"`r\")) {\n            $p.Range.Text = $r.New\n        }\n    }\n}\n"}
